$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that currently follows the
#    Heading1 title paragraph.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete() | Out-Null

# 2. Insert a new bold title paragraph right before the final paragraph
#    (which currently holds the "Create a feature image..." image prompt).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastPara.Range.InsertParagraphBefore() | Out-Null

$titlePara = $d.Paragraphs($d.Paragraphs.Count - 1)
$titleXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cleopatra II Free - Exciting Bonuses &amp; High Wins</w:t></w:r></w:p>"
$titlePara.Range.InsertXML($titleXml) | Out-Null

# 3. Replace the final paragraph's text (was the image prompt) with the
#    meta description text, keeping the italic run formatting.
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Cleopatra II offers exciting bonuses and substantial wins with high volatility and 95.13% RTP. Play for free to experience the game mechanics and Egyptian theme.</w:t></w:r></w:p>"
$finalPara.Range.InsertXML($finalXml) | Out-Null
